# TC02_Trials_Filter_Gender-Female.xlsx
# Insert a new first column (A) holding a "TabName"/"Cases" label pair,
# shifting the existing query/result columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing A:D shift to B:E.
$ws.Range("A1").EntireColumn.Insert()

# Populate the new column.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "Cases"

# Best-fit the new column's width to its (short) contents.
$ws.Columns("A").ColumnWidth = 8

# Selection ends up on B8, matching the saved view state.
$ws.Range("B8").Select()
